# Update publications_dictionary.xlsx:
# Insert a new row above the "submitted" status row (row 11) that documents
# the new "reject_resubmit" status2 value, shifting subsequent rows down by
# one, then move the active selection to the cell that replaced the old
# "submitted" row (now B12).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 11 (shifts existing rows 11-24 down to 12-25)
$ws.Rows(11).Insert(-4121)

# Copy the border/formatting from the row above (row 9, a "status2" value
# row with the same left/right-border look) onto the new row's cells
$ws.Cells.Item(9, 1).Copy()
$ws.Cells.Item(11, 1).PasteSpecial(-4122)
$ws.Cells.Item(9, 2).Copy()
$ws.Cells.Item(11, 2).PasteSpecial(-4122)

# Set the new cell's content
$ws.Cells.Item(11, 2).Value = "reject_resubmit"

# Update the selection to reflect where the user ended up after editing
$ws.Range("B12").Select() | Out-Null
